# Add a :sheet => Int option to the Excel reader spec fixture:
# duplicate "Sheet 1" as a new "Sheet 2" (same layout/styles), but with
# different sample data in the "name" column so the reader's :sheet
# option can be exercised against two distinct sheets.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet 1")

# Copy Sheet 1 right after itself - this clones all formatting, column
# widths, merged cells, page setup, etc. without having to rebuild them
# by hand.
$sheet1.Copy($null, $sheet1)

$sheet2 = $wb.Worksheets.Item(2)
$sheet2.Name = "Sheet 2"

# Give the second sheet its own sample rows.
$sheet2.Range("B3").Value = "Louis Lambeau"
$sheet2.Range("B4").Value = "Marie Deserable"

# Keep "Sheet 1" as the active/selected sheet, as in the original workbook.
$sheet1.Activate()
